# Update Work Week and Social Spending
# Refresh the "GDP per Capita" series on the Data sheet for Ireland
# (country code 372) with revised figures for 1820-2010, and extend the
# series through 2016 by appending six new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New GDP per Capita values, one per year starting at 1820 (row 2) through
# 2016 (row 198). Blank entries represent years with no recorded data.
$values = @('1398','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','2829','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','','4361','','','','','','','','4038','4141','4105','4095','4101','4100','4229','4363','4501','4618','4737','4599','4463','4594','4728','4865','4713','4865','4865','4865','4865','4865','4865','4865','4812','4865','4929','5149','5426','5504','5649','5805','5973','6048','6248','6212','6239','6169','6437','6825','7186','7390','7685','7948','8051','8097','8531','9197','9706','9881','10128','10621','10946','11225','11662','11639','12425','13150','13335','13614','13893','14061','13931','14435','14834','14768','15458','16313','17343','18838','19253.9434511805','20032.7186903767','20627.0785724041','22004.8139153003','24266.1422151915','26562.9089844149','29455.9203494766','31908.6479040328','35377.1621932825','38806.5013574153','40966.3329700636','43012.8146170694','44372.7564064811','47028.8647351758','49223.3830552297','51296.194414648','52322.2312530697','49583.1415262939','47375.7345392526','48623.8105155844','48980','48333','48743','52651','54278','56597')

$firstYear = 1820
$firstRow = 2

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $firstRow + $i
    $year = $firstYear + $i

    if ($row -gt 192) {
        # Years 2011-2016 are brand new rows; fill in the full record.
        $ws.Cells.Item($row, 1).Value = 372
        $ws.Cells.Item($row, 2).Value = "Ireland"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = $year
    }

    $ws.Cells.Item($row, 5).Value = $values[$i]
}
